$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4864.522
$ws.Range("I116").Value = 2937.5
$ws.Range("J116").Value = 5892.2666
$ws.Range("K116").Value = 2937.5
$ws.Range("L116").Value = 5892.2666
$ws.Range("M116").Value = 504.5
$ws.Range("N116").Value = -12776.2666
$ws.Range("H118").Value = 1220.1818
$ws.Range("I118").Value = 1410
$ws.Range("J118").Value = 888
$ws.Range("K118").Value = 4230
$ws.Range("L118").Value = 2664
$ws.Range("M118").Value = -2573
$ws.Range("N118").Value = -5978
$ws.Range("H132").Value = 20085.02
$ws.Range("I132").Value = 3213.7317
$ws.Range("J132").Value = 89257.3
$ws.Range("K132").Value = 9641.195099999999
$ws.Range("L132").Value = 267771.9
$ws.Range("M132").Value = -7111.195099999999
$ws.Range("N132").Value = -272831.9
$ws.Range("H135").Value = 9804736
$ws.Range("I135").Value = 733.54285
$ws.Range("K135").Value = 6601.88565
$ws.Range("M135").Value = -4066.88565
$ws.Range("H137").Value = 3700.2886
$ws.Range("I137").Value = 1270.0264
$ws.Range("J137").Value = 10296.714
$ws.Range("K137").Value = 3810.0792
$ws.Range("L137").Value = 30890.142
$ws.Range("M137").Value = -1260.0792
$ws.Range("N137").Value = -35990.142
$ws.Range("H138").Value = 1987.64
$ws.Range("I138").Value = 643.9773
$ws.Range("J138").Value = 3043.375
$ws.Range("K138").Value = 1931.9319
$ws.Range("L138").Value = 9130.125
$ws.Range("M138").Value = 3208.0681
$ws.Range("N138").Value = -19410.125
$ws.Range("H141").Value = 1188.4546
$ws.Range("I141").Value = 675.2
$ws.Range("J141").Value = 6321
$ws.Range("K141").Value = 2025.6
$ws.Range("L141").Value = 18963
$ws.Range("M141").Value = 3154.4
$ws.Range("N141").Value = -29323

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10905.639
$ws.Range("I32").Value = 11630.792
$ws.Range("J32").Value = 8882.842000000001
$ws.Range("K32").Value = 11630.792
$ws.Range("L32").Value = 8882.842000000001
$ws.Range("M32").Value = -11343.792
$ws.Range("N32").Value = -9456.842000000001
$ws.Range("H61").Value = 1115.9778
$ws.Range("I61").Value = 990.4286
$ws.Range("J61").Value = 1408.9259
$ws.Range("K61").Value = 990.4286
$ws.Range("L61").Value = 1408.9259
$ws.Range("M61").Value = -778.4286
$ws.Range("N61").Value = -1832.9259
$ws.Range("H74").Value = 1675.3684
$ws.Range("I74").Value = 1699.5272
$ws.Range("J74").Value = 1612.0952
$ws.Range("K74").Value = 1699.5272
$ws.Range("L74").Value = 1612.0952
$ws.Range("M74").Value = -825.5272
$ws.Range("N74").Value = -3360.0952
$ws.Range("H77").Value = 1675.3684
$ws.Range("I77").Value = 1699.5272
$ws.Range("J77").Value = 1612.0952
$ws.Range("K77").Value = 8497.636
$ws.Range("L77").Value = 8060.476
$ws.Range("M77").Value = -4129.636
$ws.Range("N77").Value = -16796.476
$ws.Range("H132").Value = 14708292
$ws.Range("I132").Value = 26317440
$ws.Range("K132").Value = 78952320
$ws.Range("M132").Value = -78949790
$ws.Range("H136").Value = 1115.9778
$ws.Range("I136").Value = 990.4286
$ws.Range("J136").Value = 1408.9259
$ws.Range("K136").Value = 2971.2858
$ws.Range("L136").Value = 4226.7777
$ws.Range("M136").Value = -421.2857999999997
$ws.Range("N136").Value = -9326.777699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1849.8445
$ws.Range("I105").Value = 1432.5625
$ws.Range("J105").Value = 2080.0688
$ws.Range("K105").Value = 1432.5625
$ws.Range("L105").Value = 2080.0688
$ws.Range("M105").Value = 314.4375
$ws.Range("N105").Value = -5574.0688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1780.08
$ws.Range("I31").Value = 983.0323
$ws.Range("J31").Value = 3080.5264
$ws.Range("K31").Value = 983.0323
$ws.Range("L31").Value = 3080.5264
$ws.Range("M31").Value = -688.0323
$ws.Range("N31").Value = -3670.5264
$ws.Range("H34").Value = 1780.08
$ws.Range("I34").Value = 983.0323
$ws.Range("J34").Value = 3080.5264
$ws.Range("K34").Value = 983.0323
$ws.Range("L34").Value = 3080.5264
$ws.Range("M34").Value = -781.0323
$ws.Range("N34").Value = -3484.5264
$ws.Range("H58").Value = 1276.8937
$ws.Range("I58").Value = 876.4474
$ws.Range("J58").Value = 2967.6667
$ws.Range("K58").Value = 876.4474
$ws.Range("L58").Value = 2967.6667
$ws.Range("M58").Value = -673.4474
$ws.Range("N58").Value = -3373.6667
$ws.Range("H99").Value = 1501.2222
$ws.Range("I99").Value = 1573
$ws.Range("J99").Value = 1250
$ws.Range("K99").Value = 1573
$ws.Range("L99").Value = 1250
$ws.Range("M99").Value = -75
$ws.Range("N99").Value = -4246
$ws.Range("H126").Value = 1501.2222
$ws.Range("I126").Value = 1573
$ws.Range("J126").Value = 1250
$ws.Range("K126").Value = 4719
$ws.Range("L126").Value = 3750
$ws.Range("M126").Value = -2249
$ws.Range("N126").Value = -8690
$ws.Range("H132").Value = 25137.441
$ws.Range("I132").Value = 1331.8541
$ws.Range("J132").Value = 129016.37
$ws.Range("K132").Value = 3995.5623
$ws.Range("L132").Value = 387049.11
$ws.Range("M132").Value = -1465.5623
$ws.Range("N132").Value = -392109.11
$ws.Range("H134").Value = 439478.3
$ws.Range("I134").Value = 1313.2778
$ws.Range("J134").Value = 1002833.4
$ws.Range("K134").Value = 3939.8334
$ws.Range("L134").Value = 3008500.2
$ws.Range("M134").Value = -1404.8334
$ws.Range("N134").Value = -3013570.2
$ws.Range("H136").Value = 1276.8937
$ws.Range("I136").Value = 876.4474
$ws.Range("J136").Value = 2967.6667
$ws.Range("K136").Value = 2629.3422
$ws.Range("L136").Value = 8903.000100000001
$ws.Range("M136").Value = -79.34220000000005
$ws.Range("N136").Value = -14003.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 10533.15
$ws.Range("I107").Value = 8126.769
$ws.Range("J107").Value = 15002.143
$ws.Range("K107").Value = 24380.307
$ws.Range("L107").Value = 45006.429
$ws.Range("M107").Value = -22460.307
$ws.Range("N107").Value = -48846.429
$ws.Range("H113").Value = 3138.7437
$ws.Range("I113").Value = 4542.44
$ws.Range("J113").Value = 632.1429000000001
$ws.Range("K113").Value = 13627.32
$ws.Range("L113").Value = 1896.4287
$ws.Range("M113").Value = -11457.32
$ws.Range("N113").Value = -6236.4287
$ws.Range("H122").Value = 2276.5874
$ws.Range("I122").Value = 614.8261
$ws.Range("J122").Value = 6773.1177
$ws.Range("K122").Value = 5533.4349
$ws.Range("L122").Value = 60958.0593
$ws.Range("M122").Value = -3083.4349
$ws.Range("N122").Value = -65858.05929999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 39998
$ws.Range("J93").Value = 39998
$ws.Range("L93").Value = 39998
$ws.Range("N93").Value = -43742
$ws.Range("H95").Value = 9546.4
$ws.Range("J95").Value = 9546.4
$ws.Range("L95").Value = 9546.4
$ws.Range("N95").Value = -15038.4
$ws.Range("H124").Value = 40257.332
$ws.Range("J124").Value = 40257.332
$ws.Range("L124").Value = 40257.332
$ws.Range("N124").Value = -50077.332
$ws.Range("H132").Value = 2546.5293
$ws.Range("I132").Value = 1645.2632
$ws.Range("J132").Value = 3688.1333
$ws.Range("K132").Value = 4935.7896
$ws.Range("L132").Value = 11064.3999
$ws.Range("M132").Value = -2405.7896
$ws.Range("N132").Value = -16124.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6946486.5
$ws.Range("I82").Value = 2037.625
$ws.Range("K82").Value = 2037.625
$ws.Range("M82").Value = -1676.625
$ws.Range("H85").Value = 6946486.5
$ws.Range("I85").Value = 2037.625
$ws.Range("K85").Value = 2037.625
$ws.Range("M85").Value = -789.625
$ws.Range("H132").Value = 2488.9592
$ws.Range("I132").Value = 1907.8823
$ws.Range("J132").Value = 3806.0667
$ws.Range("K132").Value = 5723.6469
$ws.Range("L132").Value = 11418.2001
$ws.Range("M132").Value = -3193.6469
$ws.Range("N132").Value = -16478.2001
$ws.Range("H136").Value = 1129.0588
$ws.Range("I136").Value = 870.36206
$ws.Range("J136").Value = 2629.5
$ws.Range("K136").Value = 2611.08618
$ws.Range("L136").Value = 7888.5
$ws.Range("M136").Value = -61.08618000000024
$ws.Range("N136").Value = -12988.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 21766
$ws.Range("J92").Value = 21766
$ws.Range("L92").Value = 21766
$ws.Range("N92").Value = -26758
$ws.Range("H132").Value = 1487.0793
$ws.Range("I132").Value = 1321.8959
$ws.Range("J132").Value = 2015.6666
$ws.Range("K132").Value = 3965.6877
$ws.Range("L132").Value = 6046.9998
$ws.Range("M132").Value = -1435.6877
$ws.Range("N132").Value = -11106.9998
$ws.Range("H136").Value = 14504
$ws.Range("I136").Value = 21391.791
$ws.Range("K136").Value = 64175.37300000001
$ws.Range("M136").Value = -61625.37300000001
